$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cells whose new numeric-looking text must stay text (force text number format
# before assignment so Excel does not auto-convert them to real numbers).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"

$ws.Range("D2").Value = '26.621.32'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '1.596.35'
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '211.16'
$ws.Range("E5").Value = '  +0.01%  '
$ws.Range("E6").Value = '  +0.96%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '0.0617'
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '0.245'
$ws.Range("E9").Value = '  -1.44%  '
$ws.Range("D10").Value = '19.38'
$ws.Range("E10").Value = '  -1.65%  '
$ws.Range("D11").Value = '0.0838'
$ws.Range("E11").Value = '  +0.33%  '
$ws.Range("D12").Value = '1.818.99'
$ws.Range("E12").Value = '  +0.23%  '
$ws.Range("D13").Value = '1.582.45'
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("E14").Value = '  -0.37%  '
$ws.Range("E15").Value = '  -1.20%  '
$ws.Range("D16").Value = '64.66'
$ws.Range("E16").Value = '  -0.20%  '
$ws.Range("D17").Value = '26.610.12'
$ws.Range("E17").Value = '  -0.16%  '
$ws.Range("E18").Value = '  +0.46%  '
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").Value = '207.79'
$ws.Range("E20").Value = '  -0.25%  '
$ws.Range("D21").Value = '6.92'
$ws.Range("E21").Value = '  +2.20%  '
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("E23").Value = '  -3.35%  '
$ws.Range("D24").Value = '8.85'
$ws.Range("E24").Value = '  -0.70%  '
$ws.Range("D25").Value = '145.47'
$ws.Range("E25").Value = '  -1.11%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").Value = '7.15'
$ws.Range("E27").Value = '  -1.78%  '
$ws.Range("E28").Value = '  +0.29%  '
$ws.Range("D29").Value = '15.29'
$ws.Range("E29").Value = '  -0.20%  '
$ws.Range("E30").Value = '  -0.19%  '
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("E32").Value = '  -0.22%  '
$ws.Range("D33").Value = '0.656'
$ws.Range("E33").Value = '  +0.52%  '
$ws.Range("E34").Value = '  +0.48%  '
$ws.Range("D35").Value = '1.282.90'
$ws.Range("E35").Value = '  -2.52%  '
$ws.Range("E36").Value = '  +1.67%  '
$ws.Range("E37").Value = '  -0.57%  '
$ws.Range("E38").Value = '  -0.34%  '
$ws.Range("D39").Value = '0.839'
$ws.Range("E39").Value = '  +1.27%  '
$ws.Range("E40").Value = '  -0.03%  '
$ws.Range("D41").Value = '5.44'
$ws.Range("E41").Value = '  +0.83%  '
$ws.Range("E42").Value = '  +1.02%  '
$ws.Range("E43").Value = '  -0.58%  '
$ws.Range("D44").Value = '63.72'
$ws.Range("E44").Value = '  +0.30%  '
$ws.Range("D45").Value = '0.919'
$ws.Range("E45").Value = '  +9.42%  '
$ws.Range("D46").Value = '1.732.00'
$ws.Range("E46").Value = '  +0.26%  '
$ws.Range("D47").Value = '89.66'
$ws.Range("E47").Value = '  -0.31%  '
$ws.Range("E48").Value = '  -0.40%  '
$ws.Range("E49").Value = '  -1.83%  '
$ws.Range("E50").Value = '  +3.46%  '
$ws.Range("E51").Value = '  -0.80%  '
